$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 85
$ws1.Range("F3").Value = 21808
$ws1.Range("F14").Value = 1033
$ws1.Range("F16").Value = 573
$ws1.Range("F23").Value = 1243
$ws1.Range("F28").Value = 624
$ws1.Range("F30").Value = 173
$ws1.Range("F31").Value = 5291
$ws1.Range("F36").Value = 13572
$ws1.Range("F38").Value = 168
$ws1.Range("F42").Value = 495
$ws1.Range("F43").Value = 4111

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 85
$ws4.Range("F3").Value = 21808
$ws4.Range("F14").Value = 1033
$ws4.Range("F16").Value = 573
$ws4.Range("F23").Value = 1243
$ws4.Range("F29").Value = 624
$ws4.Range("F32").Value = 173
$ws4.Range("F34").Value = 5291
$ws4.Range("F39").Value = 13572
$ws4.Range("F41").Value = 168
$ws4.Range("F45").Value = 495
$ws4.Range("F46").Value = 4111
